$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P4").Value = "Umid Web Site Auto_drive, [09.09.2022 21:20]"
$ws.Range("P5").Value = "3 bn 5  qala ekskursiya 45$"
$ws.Range("P6").Value = "7 qala ekskursiya - 50$"
$ws.Range("P7").Value = "10 qala - 55$"
$ws.Range("P9").Value = "Umid Web Site Auto_drive, [09.09.2022 21:22]"
$ws.Range("P10").Value = "Khiva,Buxoro, Samarqand, Toshkent  ekskursiya - 45$ 3 yil tajriba un"
$ws.Range("P11").Value = "5 yil tajriba un 49$"
$ws.Range("P12").Value = "10 yil tajriba un 59$"

$ws.Columns.Item(16).ColumnWidth = 42.0

$ws.Range("N15").Select()
